$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2025.6364
$ws.Range("I15").Value = 2025.6364
$ws.Range("K15").Value = 6076.9092
$ws.Range("M15").Value = -5907.9092
$ws.Range("H17").Value = 1971.2
$ws.Range("J17").Value = 1971.2
$ws.Range("L17").Value = 5913.6
$ws.Range("N17").Value = -6249.6
$ws.Range("H32").Value = 9096573
$ws.Range("I32").Value = 20000
$ws.Range("J32").Value = 10004230
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 10004230
$ws.Range("M32").Value = -19674
$ws.Range("N32").Value = -10004882
$ws.Range("H70").Value = 2133
$ws.Range("I70").Value = 1452.8182
$ws.Range("J70").Value = 2708.5386
$ws.Range("K70").Value = 4358.4546
$ws.Range("L70").Value = 8125.6158
$ws.Range("M70").Value = -4088.4546
$ws.Range("N70").Value = -8665.6158
$ws.Range("H73").Value = 2133
$ws.Range("I73").Value = 1452.8182
$ws.Range("J73").Value = 2708.5386
$ws.Range("K73").Value = 4358.4546
$ws.Range("L73").Value = 8125.6158
$ws.Range("M73").Value = -3422.4546
$ws.Range("N73").Value = -9997.6158
$ws.Range("H98").Value = 806.4
$ws.Range("I98").Value = 784.8889
$ws.Range("K98").Value = 784.8889
$ws.Range("M98").Value = 713.1111
$ws.Range("H99").Value = 71972
$ws.Range("I99").Value = 486.8
$ws.Range("K99").Value = 1460.4
$ws.Range("M99").Value = 37.59999999999991
$ws.Range("H100").Value = 5416.3687
$ws.Range("I100").Value = 2781.625
$ws.Range("J100").Value = 7332.5454
$ws.Range("K100").Value = 2781.625
$ws.Range("L100").Value = 7332.5454
$ws.Range("M100").Value = -2240.625
$ws.Range("N100").Value = -8414.545399999999
$ws.Range("H107").Value = 396.0909
$ws.Range("I107").Value = 358.42856
$ws.Range("K107").Value = 358.42856
$ws.Range("M107").Value = 1561.57144
$ws.Range("H112").Value = 1790.1428
$ws.Range("J112").Value = 1843.0526
$ws.Range("L112").Value = 5529.1578
$ws.Range("N112").Value = -7745.1578
$ws.Range("H122").Value = 806.4
$ws.Range("I122").Value = 784.8889
$ws.Range("K122").Value = 2354.6667
$ws.Range("M122").Value = 95.33329999999978
$ws.Range("H141").Value = 4884.0835
$ws.Range("I141").Value = 4884.0835
$ws.Range("K141").Value = 14652.2505
$ws.Range("M141").Value = -9472.250499999998

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1503.75
$ws.Range("J2").Value = 1905.5
$ws.Range("L2").Value = 1905.5
$ws.Range("N2").Value = -2131.5
$ws.Range("H116").Value = 1503.75
$ws.Range("J116").Value = 1905.5
$ws.Range("L116").Value = 1905.5
$ws.Range("N116").Value = -6493.5
$ws.Range("H132").Value = 3134
$ws.Range("I132").Value = 2480.1765
$ws.Range("K132").Value = 7440.529500000001
$ws.Range("M132").Value = -4910.529500000001

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1503.75
$ws.Range("J3").Value = 1905.5
$ws.Range("L3").Value = 1905.5
$ws.Range("N3").Value = -2133.5
$ws.Range("H94").Value = 2285
$ws.Range("I94").Value = 2102.5625
$ws.Range("J94").Value = 2702
$ws.Range("K94").Value = 2102.5625
$ws.Range("L94").Value = 2702
$ws.Range("M94").Value = -1651.5625
$ws.Range("N94").Value = -3604
$ws.Range("H122").Value = 49941.668
$ws.Range("J122").Value = 49941.668
$ws.Range("L122").Value = 49941.668
$ws.Range("N122").Value = -59741.668
$ws.Range("H134").Value = 3673.1428
$ws.Range("I134").Value = 1885.2222
$ws.Range("K134").Value = 5655.6666
$ws.Range("M134").Value = -3120.6666

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1028
$ws.Range("I16").Value = 793.1429000000001
$ws.Range("K16").Value = 793.1429000000001
$ws.Range("M16").Value = -506.1429000000001
$ws.Range("H58").Value = 2315.4167
$ws.Range("I58").Value = 1648.875
$ws.Range("K58").Value = 1648.875
$ws.Range("M58").Value = -1445.875
$ws.Range("H99").Value = 4160884.8
$ws.Range("I99").Value = 1748008.2
$ws.Range("J99").Value = 5005391.5
$ws.Range("K99").Value = 1748008.2
$ws.Range("L99").Value = 5005391.5
$ws.Range("M99").Value = -1746510.2
$ws.Range("N99").Value = -5008387.5
$ws.Range("H105").Value = 723.6667
$ws.Range("I105").Value = 646.46155
$ws.Range("K105").Value = 646.46155
$ws.Range("M105").Value = 1100.53845
$ws.Range("H113").Value = 1028
$ws.Range("I113").Value = 793.1429000000001
$ws.Range("K113").Value = 793.1429000000001
$ws.Range("M113").Value = 1376.8571
$ws.Range("H126").Value = 4160884.8
$ws.Range("I126").Value = 1748008.2
$ws.Range("J126").Value = 5005391.5
$ws.Range("K126").Value = 5244024.6
$ws.Range("L126").Value = 15016174.5
$ws.Range("M126").Value = -5241554.6
$ws.Range("N126").Value = -15021114.5
$ws.Range("H136").Value = 2315.4167
$ws.Range("I136").Value = 1648.875
$ws.Range("K136").Value = 4946.625
$ws.Range("M136").Value = -2396.625

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7692.9585
$ws.Range("I5").Value = 3731
$ws.Range("J5").Value = 10522.929
$ws.Range("K5").Value = 11193
$ws.Range("L5").Value = 31568.787
$ws.Range("M5").Value = -11081
$ws.Range("N5").Value = -31792.787
$ws.Range("H51").Value = 2000
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5540
$ws.Range("H114").Value = 1005.8
$ws.Range("I114").Value = 309.8
$ws.Range("J114").Value = 1353.8
$ws.Range("K114").Value = 929.4000000000001
$ws.Range("L114").Value = 4061.4
$ws.Range("M114").Value = 2324.6
$ws.Range("N114").Value = -10569.4
$ws.Range("H117").Value = 760.2
$ws.Range("J117").Value = 1083.75
$ws.Range("L117").Value = 3251.25
$ws.Range("N117").Value = -10135.25
$ws.Range("H132").Value = 580.7143
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H133").Value = 919
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H135").Value = 7692.9585
$ws.Range("I135").Value = 3731
$ws.Range("J135").Value = 10522.929
$ws.Range("K135").Value = 33579
$ws.Range("L135").Value = 94706.361
$ws.Range("M135").Value = -31044
$ws.Range("N135").Value = -99776.361
$ws.Range("H137").Value = 2195.2144
$ws.Range("I137").Value = 1614.7778
$ws.Range("K137").Value = 4844.3334
$ws.Range("M137").Value = 255.6665999999996
$ws.Range("N132").ClearContents()
$ws.Range("N133").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 106823.82
$ws.Range("I70").Value = 162794.58
$ws.Range("J70").Value = 8875
$ws.Range("K70").Value = 162794.58
$ws.Range("L70").Value = 8875
$ws.Range("M70").Value = -162524.58
$ws.Range("N70").Value = -9415
$ws.Range("H73").Value = 106823.82
$ws.Range("I73").Value = 162794.58
$ws.Range("J73").Value = 8875
$ws.Range("K73").Value = 162794.58
$ws.Range("L73").Value = 8875
$ws.Range("M73").Value = -161858.58
$ws.Range("N73").Value = -10747
$ws.Range("H97").Value = 381.9
$ws.Range("I97").Value = 381.9
$ws.Range("K97").Value = 381.9
$ws.Range("M97").Value = 114.1
$ws.Range("H113").Value = 4442.85
$ws.Range("I113").Value = 2714.2778
$ws.Range("K113").Value = 2714.2778
$ws.Range("M113").Value = -544.2777999999998
$ws.Range("H122").Value = 3700.3845
$ws.Range("I122").Value = 3595.7144
$ws.Range("J122").Value = 4140
$ws.Range("K122").Value = 10787.1432
$ws.Range("L122").Value = 12420
$ws.Range("M122").Value = -8337.143199999999
$ws.Range("N122").Value = -17320
$ws.Range("H135").Value = 67500
$ws.Range("J135").Value = 67500
$ws.Range("L135").Value = 67500
$ws.Range("N135").Value = -77640
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 102.6
$ws.Range("I16").Value = 102.6
$ws.Range("K16").Value = 102.6
$ws.Range("M16").Value = 67.40000000000001
$ws.Range("H45").Value = 19250
$ws.Range("I45").Value = 19250
$ws.Range("K45").Value = 19250
$ws.Range("M45").Value = -18843
$ws.Range("H46").Value = 8348.75
$ws.Range("I46").Value = 2793
$ws.Range("J46").Value = 9689.793
$ws.Range("K46").Value = 2793
$ws.Range("L46").Value = 9689.793
$ws.Range("M46").Value = -2605
$ws.Range("N46").Value = -10065.793
$ws.Range("H55").Value = 2289.0908
$ws.Range("I55").Value = 2620.111
$ws.Range("K55").Value = 2620.111
$ws.Range("M55").Value = -2447.111
$ws.Range("H93").Value = 5757.4165
$ws.Range("I93").Value = 5636.25
$ws.Range("J93").Value = 5999.75
$ws.Range("K93").Value = 5636.25
$ws.Range("L93").Value = 5999.75
$ws.Range("M93").Value = -4388.25
$ws.Range("N93").Value = -8495.75
$ws.Range("H132").Value = 3333.7659
$ws.Range("I132").Value = 2669.0588
$ws.Range("J132").Value = 5072.231
$ws.Range("K132").Value = 8007.176399999999
$ws.Range("L132").Value = 15216.693
$ws.Range("M132").Value = -5477.176399999999
$ws.Range("N132").Value = -20276.693

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 335666.34
$ws.Range("J4").Value = 3500
$ws.Range("L4").Value = 3500
$ws.Range("N4").Value = -3726
$ws.Range("H136").Value = 4783.8276
$ws.Range("I136").Value = 4774.1875
$ws.Range("J136").Value = 4795.6924
$ws.Range("K136").Value = 14322.5625
$ws.Range("L136").Value = 14387.0772
$ws.Range("M136").Value = -11772.5625
$ws.Range("N136").Value = -19487.0772
